$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Rename the worksheet (the tab name had a stray "_xlsx" suffix from the
#    original export; the author cleaned it up to match the workbook name).
# ---------------------------------------------------------------------------
$ws.Name = "class_schedule"

# ---------------------------------------------------------------------------
# 2. The GIS reading-list cell (C27) gets swapped with the Opioid-project
#    cell (C28), and a refreshed version of the GIS reading list (merging
#    the old "Vector Data" + "Raster Data" bullets into a single linked
#    "Geopandas / Vector Data" bullet) is written into C28.
#
#    We stash the original formatting (style 9: Fira Code font, wrapText,
#    quotePrefix) of C27 before overwriting it, since a plain value-write
#    resets the "quote prefix" flag that these dash-led strings rely on.
# ---------------------------------------------------------------------------

# Stash C27's current formatting (quotePrefix + wrap) off to the side.
$ws.Range("C27").Copy($ws.Range("Z1"))

# Move the Opioid-project text (with its own formatting) from C28 into C27.
$ws.Range("C28").Copy($ws.Range("C27"))

# Write the refreshed GIS reading list into C28.
$newGisText = "- ``What is GIS? <gis_what_is_gis.ipynb>``_`n" + `
    "- ``Installing Geopandas <gis_setup_geopandas.ipynb>``_`n" + `
    "- ``Geopandas / Vector Data <gis_geopandas.ipynb>``_`n" + `
    "- ``GeoPandas User Guide: Data Structures, Reading and Writing Files, and Indexing and Selecting Data <https://geopandas.org/docs/user_guide.html>``_"
$ws.Range("C28").Value = $newGisText

# Re-apply the stashed formatting (quotePrefix + wrap) onto C28, since the
# value-write above reset it, then clean up the scratch cell.
$ws.Range("Z1").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("Z1").Clear()
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Row 27 shrinks (75 -> 60) now that it holds the shorter Opioid text.
# ---------------------------------------------------------------------------
$ws.Rows(27).RowHeight = 60

# ---------------------------------------------------------------------------
# 4. Move the selection/viewport down to the newly edited row.
# ---------------------------------------------------------------------------
$ws.Range("C28").Select()
